# Atualização automática de NONOAI.xlsx
#
# Semantic changes applied:
#   1. Delete the "Desarquivamentos Pendentes" worksheet entirely.
#   2. Rename "Paineis DARQ" -> "PAINEIS DARQ".
#   3. Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO".

$wb = $excel.ActiveWorkbook

# Avoid any "are you sure you want to delete" confirmation prompts.
$excel.DisplayAlerts = $false

# 1) Remove the obsolete "Desarquivamentos Pendentes" sheet.
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete() | Out-Null

# 2) Rename "Paineis DARQ" to its all-caps form.
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"

# 3) Rename "Recolhimento x Eliminacao" to its all-caps accented form.
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $true | Out-Null
